# Apply the "Minor changes and additions" edit to panos_courses.xlsx
#
# - Changed some semester values in files-for-importing dir:
#     D2 (Year 2 course's Semester): 1 -> "Spring"
#     D3 (Year 1 course's Semester): 2 -> "Fall"
#     D4 (Year 3 course's Semester): 3 -> "Fall"
#   Also the two Syllabus texts that used to carry a stray leading space are
#   rewritten without it.
# - The defined name / query that the sheet data came from was renamed from
#   "panos_courses" to "panos_courses_1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates -------------------------------------------------
# Order matters for how new shared strings get appended: write the "Fall"
# rows before the "Spring" row, and the second syllabus text before the
# first, so new unique strings land in the same order as the target file.
$ws.Range("D3").Value = "Fall"
$ws.Range("D4").Value = "Fall"
$ws.Range("D2").Value = "Spring"

$ws.Range("F4").Value = "Advanced Databases II and more"
$ws.Range("F2").Value = "Software development basics"

# --- Rename the defined name (panos_courses -> panos_courses_1) ----------
$wb.Names.Item(1).Name = "panos_courses_1"

# --- Restore the active selection on the sheet ----------------------------
$ws.Range("E8").Select() | Out-Null
